# MAJ JDD NewAccount all YES
# Set column B (rows 3 to 9) on the "Test Cases" sheet to "YES",
# and update the active selection to B2:B9 (activeCell B2) to match
# the saved sheet view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("B3:B9").Value = "YES"

$ws.Activate()
$ws.Range("B2:B9").Select()
